# Generate Report for Handoff
# The 5da074f1-ffdd-4545-8c7d-bdfa7ab823f3.md file has moved from
# "In Translation" to "Ready for handoff" for both the zh-cn and de-de
# locales, and a new handoff datetime is recorded for each locale.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update status column for the 5da074f1 row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: update status + handoff datetime for the 5da074f1 row ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-09 20:13:16"

# --- de-de sheet: update status + handoff datetime for the 5da074f1 row ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-09 20:13:20"
